# Generate Report for Handoff
#
# The "b.md" row (row 3) on every sheet is refreshed to reflect that the
# file has now been handed off again: status becomes "Ready for handoff",
# a new handoff package (b.63290e5768f688058c7b37413b0a5c26c308f864.*) is
# recorded, and the handoff timestamps are updated.
#
# NOTE: this runtime's Hyperlink objects loaded from the source workbook
# don't expose their Address/TextToDisplay back to script, and deleting a
# single Hyperlinks.Item(...) is a no-op here - only clearing the whole
# Hyperlinks collection actually works. So for the sheets whose hyperlink
# text needs to change (zh-cn, de-de) we rebuild the full hyperlink set
# for that sheet (the underlying target addresses are unchanged, exactly
# like in the source file - only the *displayed* text for the handoff
# file link is updated).

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Sheet "Overview": row 3 corresponds to b.md
# ----------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"
$overview.Range("D3").Value = "2016-03-24 11:24:48"

# ----------------------------------------------------------------------
# Sheet "zh-cn": row 3 corresponds to b.md
# ----------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("E3").Value = "2016-03-24 11:24:43"

$zhcn.Hyperlinks.Delete()
$zhcn.Hyperlinks.Add($zhcn.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/6ff6d209de10b2743163637e9ca7a45e0973bbc3/e2e/a.md", $null, $null, "a.md") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/01687cf634aadee498224d38c146477e921f43f1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", $null, $null, "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/df36b5630cf1e6f1b4a1a9680197ce11b920dbcf/e2e/a.md", $null, $null, "a.md") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/9ab17262de4377b451268906d08252ab10c0277a/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", $null, $null, "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/6ff6d209de10b2743163637e9ca7a45e0973bbc3/e2e/b.md", $null, $null, "b.md") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/01687cf634aadee498224d38c146477e921f43f1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", $null, $null, "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/df36b5630cf1e6f1b4a1a9680197ce11b920dbcf/e2e/a.md", $null, $null, "a.md") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/9ab17262de4377b451268906d08252ab10c0277a/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", $null, $null, "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf") | Out-Null

# ----------------------------------------------------------------------
# Sheet "de-de": row 3 corresponds to b.md
# ----------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("E3").Value = "2016-03-24 11:24:48"

$dede.Hyperlinks.Delete()
$dede.Hyperlinks.Add($dede.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/6ff6d209de10b2743163637e9ca7a45e0973bbc3/e2e/a.md", $null, $null, "a.md") | Out-Null
$dede.Hyperlinks.Add($dede.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3f2d935c9c0c65a1bc1cb723e1e24b3d3c65a2bf/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", $null, $null, "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf") | Out-Null
$dede.Hyperlinks.Add($dede.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/86d01d4474f3a3e14ba63e25d272e9b3acaedf36/e2e/a.md", $null, $null, "a.md") | Out-Null
$dede.Hyperlinks.Add($dede.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/3a1ecaedd9fe0fd8e0ad9ef65b24535e24354340/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", $null, $null, "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf") | Out-Null
$dede.Hyperlinks.Add($dede.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/6ff6d209de10b2743163637e9ca7a45e0973bbc3/e2e/b.md", $null, $null, "b.md") | Out-Null
$dede.Hyperlinks.Add($dede.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3f2d935c9c0c65a1bc1cb723e1e24b3d3c65a2bf/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", $null, $null, "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf") | Out-Null
$dede.Hyperlinks.Add($dede.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/86d01d4474f3a3e14ba63e25d272e9b3acaedf36/e2e/a.md", $null, $null, "a.md") | Out-Null
$dede.Hyperlinks.Add($dede.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/3a1ecaedd9fe0fd8e0ad9ef65b24535e24354340/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", $null, $null, "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf") | Out-Null

Write-Output "Handoff report regenerated for b.md"
